$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: J1 connector revised from USB2.0 Micro-B to USB2.0 Type-C ---
$ws.Range("B2").Value = "USB4105-GF-A"
$ws.Range("D2").Value = "CONN RCP USB2.0 TYP C 24P SMD RA"

# --- New row 16: R8,R9 - 5.1K resistors added to BOM ---
# Copy formatting from the row above (row 15) so the new row matches the
# existing table styling (borders / alignment / fonts), then fill in values.
$ws.Range("A15:D15").Copy()
$ws.Range("A16:D16").PasteSpecial(-4122)
$ws.Range("A16").Value = "R8,R9"
$ws.Range("C16").Value = 2
$ws.Range("B16").Value = "ERA-2APB512X"
$ws.Range("D16").Value = "RES SMD 5.1K OHM 0.1% 1/16W 0402"

# --- Defined name now covers the extra row ---
$wb.Names.Item(1).RefersTo = '=Sheet1!$A$1:$D$16'

# --- Selection moved back to the top of the table ---
[void]$ws.Range("A2").Select()

# --- Page setup touched (portrait orientation stamped onto the sheet) ---
$ws.PageSetup.Orientation = 1
